$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O9").Value = 0.001238584518432617
$ws.Range("O10").Value = 0.0009899139404296875
$ws.Range("O11").Value = 0.02717375755310059
$ws.Range("O12").Value = 0.002861261367797852
$ws.Range("O13").Value = 0.01743292808532715
$ws.Range("O14").Value = 0.123436450958252
$ws.Range("O15").Value = 0.003000974655151367
$ws.Range("O16").Value = 0.006684064865112305
$ws.Range("O17").Value = 0.04587340354919434
$ws.Range("O18").Value = 0.002004146575927734
$ws.Range("O19").Value = 0.02090883255004883
$ws.Range("O20").Value = 0.0009725093841552734
$ws.Range("O21").Value = 0.01871132850646973
$ws.Range("O22").Value = 0.001993656158447266
$ws.Range("O23").Value = 0
$ws.Range("O24").Value = 0.0009891986846923828
$ws.Range("O25").Value = 0.007002830505371094
$ws.Range("O26").Value = 0.06189751625061035
$ws.Range("O27").Value = 0.002984762191772461
$ws.Range("O30").Value = 0.0030059814453125
$ws.Range("O31").Value = 0.0363919734954834
$ws.Range("O33").Value = 0.0009946823120117188
$ws.Range("O34").Value = 0.05583310127258301
$ws.Range("O35").Value = 0.02160072326660156
$ws.Range("O36").Value = 0.0008256435394287109
$ws.Range("O37").Value = 0.08654332160949707
$ws.Range("O38").Value = 0.04495835304260254
$ws.Range("O39").Value = 0.06683254241943359
$ws.Range("O40").Value = 0.04353547096252441
$ws.Range("O41").Value = 0.01604652404785156
$ws.Range("O42").Value = 0.02392148971557617
$ws.Range("O43").Value = 0.0176842212677002
$ws.Range("O44").Value = 0.03731489181518555
$ws.Range("O45").Value = 0.002000331878662109
$ws.Range("O46").Value = 0.03380107879638672
$ws.Range("O47").Value = 0.002107381820678711
$ws.Range("O48").Value = 0.01739501953125
$ws.Range("O49").Value = 0.01351213455200195
$ws.Range("O50").Value = 0.00199127197265625
$ws.Range("O51").Value = 0.006040573120117188
$ws.Range("O52").Value = 0.002001285552978516
$ws.Range("O53").Value = 0.01565098762512207
$ws.Range("O54").Value = 0.002007961273193359
$ws.Range("O55").Value = 0.04849100112915039
$ws.Range("O57").Value = 0.02603363990783691
$ws.Range("O58").Value = 0.01462340354919434
$ws.Range("O60").Value = 0.0492253303527832
$ws.Range("O62").Value = 0.006050348281860352
$ws.Range("O63").Value = 0.07803916931152344
$ws.Range("O64").Value = 0.06633329391479492
$ws.Range("O65").Value = 0.4735481739044189
$ws.Range("O66").Value = 0.01341080665588379
$ws.Range("O67").Value = 0.003023862838745117
$ws.Range("O68").Value = 0.001993656158447266
$ws.Range("O69").Value = 0.05604648590087891
$ws.Range("O70").Value = 0.04452061653137207
$ws.Range("O71").Value = 16.38140368461609
$ws.Range("O72").Value = 186.1293351650238
$ws.Range("O73").Value = 0
$ws.Range("O74").Value = 0.001004219055175781
$ws.Range("O75").Value = 0.002202749252319336
$ws.Range("O77").Value = 0
$ws.Range("O79").Value = 0.002198457717895508
$ws.Range("O80").Value = 0.00198674201965332
$ws.Range("O81").Value = 0.0674440860748291
$ws.Range("O82").Value = 0
$ws.Range("O83").Value = 0.002001762390136719
$ws.Range("O84").Value = 0.5386946201324463
$ws.Range("O85").Value = 0.00502324104309082
$ws.Range("O86").Value = 0.00593876838684082
$ws.Range("O87").Value = 0.1505951881408691
$ws.Range("O88").Value = 0.1274991035461426
$ws.Range("O89").Value = 0.004006147384643555
$ws.Range("O90").Value = 0.419173002243042
$ws.Range("O91").Value = 0.06857752799987793
$ws.Range("O92").Value = 0.06133031845092773
$ws.Range("O93").Value = 0.01563596725463867
$ws.Range("O94").Value = 0.003249406814575195
$ws.Range("O95").Value = 0.0103905200958252
$ws.Range("O96").Value = 2.70869255065918
$ws.Range("O97").Value = 0.0217583179473877
$ws.Range("O98").Value = 0.06321144104003906
$ws.Range("O99").Value = 0.05058979988098145
$ws.Range("O100").Value = 0.0568697452545166
$ws.Range("O101").Value = 0.03600049018859863
$ws.Range("O103").Value = 0.03467202186584473
$ws.Range("O104").Value = 2.23248291015625
$ws.Range("O108").Value = 0.04092812538146973
$ws.Range("O109").Value = 0.003993511199951172
$ws.Range("O110").Value = 0.05959296226501465
$ws.Range("O111").Value = 0.05014777183532715
$ws.Range("O112").Value = 0.04639625549316406
$ws.Range("O113").Value = 0.0009820461273193359
$ws.Range("O114").Value = 0.04128050804138184
$ws.Range("O115").Value = 0.02747178077697754
$ws.Range("O116").Value = 0.0009164810180664062
$ws.Range("O117").Value = 0.002113819122314453
$ws.Range("O118").Value = 0.0009980201721191406
$ws.Range("O119").Value = 0
$ws.Range("O120").Value = 0.001997232437133789
$ws.Range("O121").Value = 0.02913784980773926
$ws.Range("O122").Value = 0.03264355659484863
$ws.Range("O124").Value = 0.00244140625
$ws.Range("O125").Value = 0.006496191024780273
$ws.Range("O126").Value = 0
$ws.Range("O127").Value = 0.001000642776489258
$ws.Range("O129").Value = 0.001044034957885742
$ws.Range("O130").Value = 0.002027511596679688
$ws.Range("O131").Value = 0.0009920597076416016
$ws.Range("O132").Value = 0.009003639221191406
$ws.Range("O133").Value = 0.001764297485351562
$ws.Range("O134").Value = 0.00101017951965332
$ws.Range("O135").Value = 0.002320766448974609
$ws.Range("O137").Value = 0
$ws.Range("O139").Value = 0.01157760620117188
$ws.Range("O140").Value = 0.03242611885070801
$ws.Range("O141").Value = 0.002999067306518555
$ws.Range("O142").Value = 0.05230903625488281
$ws.Range("O143").Value = 0.001664876937866211
$ws.Range("O144").Value = 0.0202934741973877
$ws.Range("O145").Value = 0.003092765808105469
$ws.Range("O146").Value = 0.002726316452026367
$ws.Range("O147").Value = 0.1100273132324219
$ws.Range("O148").Value = 0.03312110900878906
$ws.Range("O149").Value = 0.2554035186767578
$ws.Range("O150").Value = 0.2010571956634521
$ws.Range("O151").Value = 4.446112155914307
$ws.Range("O152").Value = 0.106736421585083
$ws.Range("O153").Value = 0
$ws.Range("O154").Value = 0.004574060440063477
$ws.Range("O155").Value = 0.1401519775390625
$ws.Range("O156").Value = 0
$ws.Range("O157").Value = 0.01747679710388184
$ws.Range("O158").Value = 0.1935446262359619
$ws.Range("O159").Value = 0.1466183662414551
$ws.Range("O160").Value = 0.004705429077148438
$ws.Range("O162").Value = 0.01831221580505371
$ws.Range("O163").Value = 0.007512092590332031
$ws.Range("O164").Value = 0.1833932399749756
$ws.Range("O165").Value = 0.120621919631958
$ws.Range("O166").Value = 0.0722966194152832
$ws.Range("O167").Value = 0.1333177089691162
$ws.Range("O169").Value = 0.01279473304748535
$ws.Range("O170").Value = 0
$ws.Range("O171").Value = 0.02595090866088867
$ws.Range("O172").Value = 0.02843761444091797
$ws.Range("O173").Value = 0.04254531860351562
$ws.Range("O174").Value = 0.05129647254943848
$ws.Range("O175").Value = 0.001169204711914062
$ws.Range("O176").Value = 0.01341414451599121
$ws.Range("O177").Value = 0.003011703491210938
